# Week 17 data logging for 2021 Team Data workbook.
# Appends this week's per-play (YDS) and per-game (ST) log entries, and
# updates the season-cumulative totals on OFF / DEF / ST / TURNS / PEN.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append this week's play-by-play yardage log to each of the
# four running logs (Rush-Off, Pass-Off, Rush-Def, Pass-Def).
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Text + " 55 2 0 1 -1 14 2 1 4 12 1 9 0 -1 0 1 0 0 11 7 2 4 22 -1 5"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Text + " 9 9 13 5 9 12 5 9 24 11 24 19 2 23 11 21 14 4"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Text + " 6 3 4 6 6 -1 7 5 4 -5 5 -2 5 -4 1 2 3 11 4 0 2"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Text + " 1 4 9 15 7 4 31 21 3 2 21 14 6 24 14 9 8 20 14 1 8 32 9 21 3 12 8 3 6 6 27 10 33"

# ---------------------------------------------------------------------
# OFF sheet: season cumulative offensive down/distance + play totals.
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("B2").Value = 3
$offWs.Range("C2").Value = 181
$offWs.Range("D2").Value = 12
$offWs.Range("E2").Value = 2
$offWs.Range("F2").Value = 40
$offWs.Range("G2").Value = 64
$offWs.Range("J2").Value = 24
$offWs.Range("L2").Value = 333
$offWs.Range("M2").Value = 201
$offWs.Range("O2").Value = 27
$offWs.Range("Q2").Value = 568

$offWs.Range("C3").Value = 191
$offWs.Range("D3").Value = 2
$offWs.Range("E3").Value = 37
$offWs.Range("F3").Value = 121
$offWs.Range("I3").Value = 70
$offWs.Range("J3").Value = 55
$offWs.Range("N3").Value = 18

# ---------------------------------------------------------------------
# DEF sheet: season cumulative defensive down/distance + play totals.
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("B2").Value = 4
$defWs.Range("C2").Value = 238
$defWs.Range("F2").Value = 67
$defWs.Range("G2").Value = 63
$defWs.Range("J2").Value = 33
$defWs.Range("L2").Value = 311
$defWs.Range("M2").Value = 219
$defWs.Range("O2").Value = 20
$defWs.Range("P2").Value = 13
$defWs.Range("Q2").Value = 607

$defWs.Range("C3").Value = 159
$defWs.Range("E3").Value = 32
$defWs.Range("F3").Value = 107
$defWs.Range("G3").Value = 35
$defWs.Range("H3").Value = 23
$defWs.Range("I3").Value = 60
$defWs.Range("J3").Value = 58
$defWs.Range("N3").Value = 23

# ---------------------------------------------------------------------
# ST sheet: season cumulative special-teams totals plus per-game logs
# for kickoff touchbacks/downed/returned/return-yards.
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 64
$stWs.Range("D2").Value = 59
$stWs.Range("F2").Value = 25
$stWs.Range("G2").Value = 21
$stWs.Range("N2").Value = 4
$stWs.Range("O2").Value = 1
$stWs.Range("B3").Value = 44

$stWs.Range("B4").Value = $stWs.Range("B4").Text + " 62"
$stWs.Range("B5").Value = $stWs.Range("B5").Text + " 7"
$stWs.Range("B6").Value = $stWs.Range("B6").Text + " 26 27 13 25 28"
$stWs.Range("D3").Value = $stWs.Range("D3").Text + " 38 32 46 38"
$stWs.Range("D4").Value = $stWs.Range("D4").Text + " 0 0 3 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Text + " 0 0 0"

# ---------------------------------------------------------------------
# TURNS sheet: season cumulative takeaway/giveaway totals.
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("C2").Value = 4
$turnsWs.Range("E3").Value = 11

# ---------------------------------------------------------------------
# PEN sheet: season cumulative penalty totals.
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 11
